# Apply edits described by the diff to the "MLK_PMT_10103_-_V-003" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$deg = [char]0x00B0

# Row 8 : clear the "Expansion Tank" description (D8:D12 is merged, keep style)
$ws.Range("D8").Value = ""

# Row 8
$ws.Range("F8").Value = "Liquid"
$ws.Range("G8").Value = "Condensate"
$ws.Range("K8").Value = "No"
$ws.Range("L8").Value = "200$deg C"
$ws.Range("N8").Value = "185$deg C"

# Row 9
$ws.Range("G9").Value = "Condensate"
$ws.Range("J9").Value = "10"
$ws.Range("K9").Value = "No"
$ws.Range("L9").Value = "200$deg C"
$ws.Range("N9").Value = "185$deg C"

# Row 10
$ws.Range("G10").Value = "Condensate"
$ws.Range("J10").Value = "Gr.B"
$ws.Range("K10").Value = "No"
$ws.Range("L10").Value = "200$deg C"
$ws.Range("N10").Value = "185$deg C"

# Row 11
$ws.Range("G11").Value = "Condensate"
$ws.Range("J11").Value = "23"
$ws.Range("K11").Value = "No"
$ws.Range("L11").Value = "200$deg C"
$ws.Range("N11").Value = "185$deg C"

# Row 12
$ws.Range("G12").Value = "Condensate"
$ws.Range("K12").Value = "No"
$ws.Range("L12").Value = "200$deg C"
$ws.Range("N12").Value = "185$deg C"
